$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Diebold-Mariano statistic (C) and P-Value (D) figures per
# "Correccion a Diebold Mariano y revision de Cap1"

$ws.Range("C2").Value = 0.9592697699939519
$ws.Range("D2").Value = 0.3441951932868075

$ws.Range("C3").Value = 0.4826817026439854
$ws.Range("D3").Value = 0.6324151094773538

$ws.Range("C4").Value = 0.423929700445918
$ws.Range("D4").Value = 0.674287262009214

$ws.Range("C5").Value = 0.1898368015468584
$ws.Range("D5").Value = 0.8505658598184254

$ws.Range("C6").Value = -0.6610682953822622
$ws.Range("D6").Value = 0.5130234776394342

$ws.Range("C7").Value = -0.7738502890056013
$ws.Range("D7").Value = 0.4443670767931187

$ws.Range("C8").Value = -0.6262124666106097
$ws.Range("D8").Value = 0.535357439098487

$ws.Range("C9").Value = -0.05157100410062943
$ws.Range("D9").Value = 0.9591723258621734

$ws.Range("C10").Value = -0.2267297295764676
$ws.Range("D10").Value = 0.8219915974713454

$ws.Range("C11").Value = -0.05723441382199545
$ws.Range("D11").Value = 0.9546935120213069
